$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = "System, dnasr281@gmail.com"
$ws.Range("G3").Value = "System, dnasr281@gmail.com"
$ws.Range("G4").Value = "System, dnasr281@gmail.com"
$ws.Range("G5").Value = "System, dnasr281@gmail.com"
$ws.Range("G6").Value = "System, dnasr281@gmail.com"
$ws.Range("G7").Value = "System, dnasr281@gmail.com"
$ws.Range("G16").Value = "System, dnasr281@gmail.com"
$ws.Range("G17").Value = "System, dnasr281@gmail.com"
$ws.Range("G22").Value = "System, dnasr281@gmail.com"
$ws.Range("G23").Value = "System, dnasr281@gmail.com"
$ws.Range("G37").Value = "System, dnasr281@gmail.com"
$ws.Range("G38").Value = "System, dnasr281@gmail.com"
$ws.Range("G43").Value = "System, dnasr281@gmail.com"
$ws.Range("G44").Value = "System, dnasr281@gmail.com"
$ws.Range("G58").Value = "System, dnasr281@gmail.com"
$ws.Range("G59").Value = "System, dnasr281@gmail.com"
$ws.Range("G64").Value = "System, dnasr281@gmail.com"
$ws.Range("G65").Value = "System, dnasr281@gmail.com"
$ws.Range("G79").Value = "System, dnasr281@gmail.com"
$ws.Range("G80").Value = "System, dnasr281@gmail.com"
$ws.Range("G85").Value = "System, dnasr281@gmail.com"
$ws.Range("G86").Value = "System, dnasr281@gmail.com"
$ws.Range("G87").Value = "System, dnasr281@gmail.com"
$ws.Range("G88").Value = "System, dnasr281@gmail.com"
$ws.Range("G89").Value = "System, dnasr281@gmail.com"
$ws.Range("G90").Value = "System, dnasr281@gmail.com"
$ws.Range("G99").Value = "System, dnasr281@gmail.com"
$ws.Range("G100").Value = "System, dnasr281@gmail.com"
$ws.Range("G105").Value = "System, dnasr281@gmail.com"
$ws.Range("G106").Value = "System, dnasr281@gmail.com"
$ws.Range("G107").Value = "System, dnasr281@gmail.com"
$ws.Range("G108").Value = "System, dnasr281@gmail.com"
$ws.Range("G109").Value = "System, dnasr281@gmail.com"
$ws.Range("G110").Value = "System, dnasr281@gmail.com"
$ws.Range("G119").Value = "System, dnasr281@gmail.com"
$ws.Range("G120").Value = "System, dnasr281@gmail.com"
$ws.Range("G125").Value = "System, dnasr281@gmail.com"
$ws.Range("G126").Value = "System, dnasr281@gmail.com"
$ws.Range("G127").Value = "System, dnasr281@gmail.com"
$ws.Range("G128").Value = "System, dnasr281@gmail.com"
$ws.Range("G129").Value = "System, dnasr281@gmail.com"
$ws.Range("G130").Value = "System, dnasr281@gmail.com"
$ws.Range("G139").Value = "System, dnasr281@gmail.com"
$ws.Range("G140").Value = "System, dnasr281@gmail.com"
$ws.Range("G145").Value = "System, dnasr281@gmail.com"
$ws.Range("G146").Value = "System, dnasr281@gmail.com"
$ws.Range("G147").Value = "System, dnasr281@gmail.com"
$ws.Range("G148").Value = "System, dnasr281@gmail.com"
$ws.Range("G149").Value = "System, dnasr281@gmail.com"
$ws.Range("G150").Value = "System, dnasr281@gmail.com"
$ws.Range("G159").Value = "System, dnasr281@gmail.com"
$ws.Range("G160").Value = "System, dnasr281@gmail.com"
$ws.Range("G165").Value = "System, dnasr281@gmail.com"
$ws.Range("G166").Value = "System, dnasr281@gmail.com"
$ws.Range("G167").Value = "System, dnasr281@gmail.com"
$ws.Range("G168").Value = "System, dnasr281@gmail.com"
$ws.Range("G169").Value = "System, dnasr281@gmail.com"
$ws.Range("G170").Value = "System, dnasr281@gmail.com"
$ws.Range("G179").Value = "System, dnasr281@gmail.com"
$ws.Range("G180").Value = "System, dnasr281@gmail.com"
$ws.Range("G185").Value = "System, dnasr281@gmail.com"
$ws.Range("G186").Value = "System, dnasr281@gmail.com"
$ws.Range("G200").Value = "System, dnasr281@gmail.com"
$ws.Range("G201").Value = "System, dnasr281@gmail.com"
$ws.Range("G206").Value = "System, dnasr281@gmail.com"
$ws.Range("G207").Value = "System, dnasr281@gmail.com"
$ws.Range("G221").Value = "System, dnasr281@gmail.com"
$ws.Range("G222").Value = "System, dnasr281@gmail.com"
$ws.Range("G227").Value = "System, dnasr281@gmail.com"
$ws.Range("G228").Value = "System, dnasr281@gmail.com"
$ws.Range("G242").Value = "System, dnasr281@gmail.com"
$ws.Range("G243").Value = "System, dnasr281@gmail.com"
